$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8, shifting the existing
# rows 8 and 9 (and their formatting, incl. the date style on column D)
# down to rows 9 and 10.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with the new weekly price record
# (week of 2022-11-03 / serial 44868).
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44868
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100114007
$ws.Range("G8").Value = "Jengibre"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 18000
$ws.Range("N8").Value = "$/caja 13 kilos"
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 1385
$ws.Range("Q8").Value = 13
$ws.Range("R8").Value = "Hortaliza"
